$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows above the old row 9 ("В0318 ...") by copying the
#    format of that row, so the existing rows 9-19 shift down to rows
#    11-21 and the two fresh rows (9-10) inherit correct per-column styles.
$ws.Rows("9:9").Copy()
$ws.Rows("9:10").Insert()
$excel.CutCopyMode = $false

# 2) Append three new rows at the end of the table (rows 22-24), copying
#    the format of the current last row (21).
$ws.Rows("21:21").Copy()
$ws.Rows("22:24").Insert()
$excel.CutCopyMode = $false

# 3) Replace old placeholder text "<нет ФН>" with "(нет ФН)" for rows 2-8 (col B)
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "(нет ФН)"
}

# 4) Write the new last row ("Д0611 ... / ГиД") first, at its final position (row 22)
$ws.Cells.Item(22, 1).Value = 'Д0611 - КИС "Единая ГИС"'
$ws.Cells.Item(22, 2).Value = "ГиД"

# 5) Write the two new "RCBI" rows. Column A of row 10 is written before
#    column A of row 9, then column B is written (shared by both rows).
$ws.Cells.Item(10, 1).Value = "А0602 - Адм. Время. RCBI. Отпуск"
$ws.Cells.Item(9, 1).Value = "А0603 - Адм. Время. RCBI. Больничный"
$ws.Cells.Item(9, 2).Value = "MES"
$ws.Cells.Item(10, 2).Value = "MES"

# 6) Write the remaining two new rows at the bottom of the table
$ws.Cells.Item(23, 1).Value = "П0558 - Отраслевой шаблон ГиД на S/4HANA"
$ws.Cells.Item(23, 2).Value = "ERP"

$ws.Cells.Item(24, 1).Value = "Т0598 - Система учета трудозатрат ЛТ"
$ws.Cells.Item(24, 2).Value = "КИС"

$ws.Range("B24").Select()
